$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.15,   4350.22, 131.82, 4797.24),
    @(0.15,   4749.3,  143.92, 542.73),
    @(0.15,   399.08,  12.1,   -4254.51),
    @(0.1,    4513.4,  136.77, 4990.83),
    @(0.1,    985.39,  29.86,  1018.67),
    @(0.1,    -3528.01,-106.91,-3972.16),
    @(0.01,   4764.56, 144.38, 4692.04),
    @(0.01,   2736.93, 82.94,  2967.03),
    @(0.01,   -2027.63,-61.44, -1725.01),
    @(0.001,  1300.78, 39.42,  2735.17),
    @(0.001,  2662.31, 80.68000000000001, 1689.51),
    @(0.001,  1361.53, 41.26,  -1045.66),
    @(0.0001, 724.1900000000001, 21.95, 920.91),
    @(0.0001, 864.77,  26.21,  784.4400000000001),
    @(0.0001, 140.58,  4.26,   -136.47)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
